# RevlonUS Latest Changes and Expansion
# Adds three new PaymentDetails rows (Discover, AMEX, MasterCard) to the DataSet sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: Discover card ---
$ws.Range("A21").Value = "PaymentDetailsDiscoverCard"
$ws.Range("N21").Value = "Discover"
$ws.Range("Q21").Value = "Mar"

# --- Row 22: AMEX card ---
$ws.Range("A22").Value = "PaymentDetailsAMEXCard"
$ws.Range("N22").Value = "AMEX"
$ws.Range("Q22").Value = "Mar"

# --- Row 23: MasterCard card ---
$ws.Range("A23").Value = "PaymentDetailsMasterCard"
$ws.Range("N23").Value = "MasterCard"
$ws.Range("Q23").Value = "Mar"

# Card numbers entered last, in MasterCard / Discover / AMEX order
# (typed with a leading apostrophe so the long digit strings are stored as text)
$ws.Range("O23").Value = "'5500005555555559"
$ws.Range("O21").Value = "'6011000995500000"
$ws.Range("O22").Value = "'371449635398431"

# Expiration years
$ws.Range("P21").Value = 2030
$ws.Range("P22").Value = 2030
$ws.Range("P23").Value = 2030

# CVVs
$ws.Range("R21").Value = 737
$ws.Range("R22").Value = 7371
$ws.Range("R23").Value = 737

# Touch the trailing date-tracking columns so the row extends through column AA,
# matching the layout of the other data rows, then leave them blank.
$ws.Range("Y21:AA21").Style = "Normal"
$ws.Range("Y22:AA22").Style = "Normal"
$ws.Range("Y23:AA23").Style = "Normal"

# Columns A and O grew / shrank to fit the newly entered values
$ws.Columns.Item(1).ColumnWidth = 27.140625
$ws.Columns.Item(15).ColumnWidth = 17.28515625

# Scroll/zoom/selection as left by the editor
$ws.Activate() | Out-Null
$ws.Range("S25").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
